$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'318.98"
$ws.Range("E2").Value = "'5.19%"
$ws.Range("G2").Value = "'15"
$ws.Range("D3").Value = "'48.81"
$ws.Range("E3").Value = "'13.46%"
$ws.Range("G3").Value = "'15"
$ws.Range("D4").Value = "'5.268"
$ws.Range("E4").Value = "'3.55%"
$ws.Range("G4").Value = "'15"
$ws.Range("D5").Value = "'0.07925"
$ws.Range("E5").Value = "'3.38%"
$ws.Range("G5").Value = "'15"
$ws.Range("D6").Value = "'4.581"
$ws.Range("E6").Value = "'3.60%"
$ws.Range("G6").Value = "'15"
$ws.Range("D7").Value = "'1.359"
$ws.Range("E7").Value = "'32.94%"
$ws.Range("G7").Value = "'15"
$ws.Range("D8").Value = "'1.637"
$ws.Range("E8").Value = "'1.47%"
$ws.Range("G8").Value = "'15"
$ws.Range("D9").Value = "'0.1289"
$ws.Range("E9").Value = "'3.81%"
$ws.Range("G9").Value = "'15"
$ws.Range("D10").Value = "'0.1958"
$ws.Range("E10").Value = "'4.82%"
$ws.Range("G10").Value = "'15"
$ws.Range("D11").Value = "'0.09454"
$ws.Range("E11").Value = "'3.21%"
$ws.Range("G11").Value = "'15"
$ws.Range("E12").Value = "'10.42%"
$ws.Range("G12").Value = "'15"
$ws.Range("D13").Value = "'0.1047"
$ws.Range("E13").Value = "'0.07%"
$ws.Range("G13").Value = "'15"
$ws.Range("D14").Value = "'0.001321"
$ws.Range("E14").Value = "'3.05%"
$ws.Range("G14").Value = "'15"
$ws.Range("D15").Value = "'0.04172"
$ws.Range("E15").Value = "'0.20%"
$ws.Range("G15").Value = "'15"
$ws.Range("D16").Value = "'0.005877"
$ws.Range("E16").Value = "'2.12%"
$ws.Range("G16").Value = "'15"
$ws.Range("B17").Value = "'HotbitToken"
$ws.Range("C17").Value = "'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D17").Value = "'0.004257"
$ws.Range("E17").Value = "'-5.18%"
$ws.Range("G17").Value = "'15"
$ws.Range("B18").Value = "'LEO"
$ws.Range("C18").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.341"
$ws.Range("E18").Value = "'0.27%"
$ws.Range("G18").Value = "'15"
$ws.Range("B19").Value = "'BTSEToken"
$ws.Range("C19").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.421"
$ws.Range("E19").Value = "'2.80%"
$ws.Range("G19").Value = "'15"
$ws.Range("B20").Value = "'BitpandaEcosystemToken"
$ws.Range("C20").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3458"
$ws.Range("E20").Value = "'3.22%"
$ws.Range("G20").Value = "'15"
$ws.Range("B21").Value = "'MCDex"
$ws.Range("C21").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "'8.190"
$ws.Range("E21").Value = "'-4.28%"
$ws.Range("G21").Value = "'15"
$ws.Range("B22").Value = "'ProBitToken"
$ws.Range("C22").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").Value = "'0.1382"
$ws.Range("E22").Value = "'-1.21%"
$ws.Range("G22").Value = "'15"
$ws.Range("B23").Value = "'ZBToken"
$ws.Range("C23").Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.3095"
$ws.Range("E23").Value = "'-3.22%"
$ws.Range("G23").Value = "'15"
$ws.Range("B24").Value = "'BitKan"
$ws.Range("C24").Value = "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "'0.001319"
$ws.Range("E24").Value = "'2.72%"
$ws.Range("G24").Value = "'15"
$ws.Range("D25").Value = "'0.0001350"
$ws.Range("E25").Value = "'-0.03%"
$ws.Range("G25").Value = "'15"
$ws.Range("D26").Value = "'0.0003541"
$ws.Range("E26").Value = "'-95.24%"
$ws.Range("G26").Value = "'15"
$ws.Range("G27").Value = "'15"
$ws.Range("G28").Value = "'15"
$ws.Range("G29").Value = "'15"
$ws.Range("G30").Value = "'15"
$ws.Range("G31").Value = "'15"
$ws.Range("G32").Value = "'15"
$ws.Range("G33").Value = "'15"
$ws.Range("G34").Value = "'15"
$ws.Range("G35").Value = "'15"
$ws.Range("G36").Value = "'15"
$ws.Range("G37").Value = "'15"
$ws.Range("D38").Value = "'0.02668"
$ws.Range("E38").Value = "'9.13%"
$ws.Range("G38").Value = "'15"
$ws.Range("D39").Value = "'0.05773"
$ws.Range("E39").Value = "'9.52%"
$ws.Range("G39").Value = "'15"
$ws.Range("D40").Value = "'0.01073"
$ws.Range("E40").Value = "'79.91%"
$ws.Range("G40").Value = "'15"
$ws.Range("D41").Value = "'0.008011"
$ws.Range("E41").Value = "'4.24%"
$ws.Range("G41").Value = "'15"
$ws.Range("D42").Value = "'0.1436"
$ws.Range("E42").Value = "'6.82%"
$ws.Range("G42").Value = "'15"
$ws.Range("D43").Value = "'0.007646"
$ws.Range("E43").Value = "'3.75%"
$ws.Range("G43").Value = "'15"
$ws.Range("D44").Value = "'0.008467"
$ws.Range("E44").Value = "'11.68%"
$ws.Range("G44").Value = "'15"
$ws.Range("D45").Value = "'0.3182"
$ws.Range("E45").Value = "'5.33%"
$ws.Range("G45").Value = "'15"
$ws.Range("D46").Value = "'0.00006631"
$ws.Range("E46").Value = "'-1.32%"
$ws.Range("G46").Value = "'15"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.02%"
$ws.Range("G47").Value = "'15"
$ws.Range("D48").Value = "'0.05492"
$ws.Range("E48").Value = "'22.44%"
$ws.Range("G48").Value = "'15"
$ws.Range("D49").Value = "'0.004002"
$ws.Range("E49").Value = "'-4.74%"
$ws.Range("G49").Value = "'15"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'0.02%"
$ws.Range("G50").Value = "'15"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'0.02%"
$ws.Range("G51").Value = "'15"
